# Add cantrals by cantons
# The sheet originally had two header rows (row 1: partial units, row 2: more
# units) above the data rows (3..10). The edit turns row 2 into the single
# header row (adding idx/idx2/Name/Date Start/Date End columns and renaming
# the power/energy headers) and removes the old row 1, shifting the data
# rows up by one (3..10 -> 2..9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the old first header row. This shifts row 2 (m3/s / MW / GWh units)
# up to row 1, and the data rows 3:10 up to 2:9.
$ws.Rows.Item(1).Delete()

# --- New row 1 header ---------------------------------------------------
# Columns A:E are brand new (no special style).
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

# Column F keeps its existing text "(m3/s)" but needs the new style (same
# font as the other header cells, default/general number format).
$ws.Range("F1").Font.Name = "Arial"
$ws.Range("F1").Font.Size = 9
$ws.Range("F1").IndentLevel = 0

# Columns G:K get new, more descriptive header text plus the same new style.
$ws.Range("G1").Value = "(MW1)"
$ws.Range("G1").Font.Name = "Arial"
$ws.Range("G1").Font.Size = 9
$ws.Range("G1").IndentLevel = 0

$ws.Range("H1").Value = "(MW2)"
$ws.Range("H1").Font.Name = "Arial"
$ws.Range("H1").Font.Size = 9
$ws.Range("H1").IndentLevel = 0

$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("I1").Font.Name = "Arial"
$ws.Range("I1").Font.Size = 9
$ws.Range("I1").IndentLevel = 0

$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("J1").Font.Name = "Arial"
$ws.Range("J1").Font.Size = 9
$ws.Range("J1").IndentLevel = 0

$ws.Range("K1").Value = "(GWh) Year"
$ws.Range("K1").Font.Name = "Arial"
$ws.Range("K1").Font.Size = 9
$ws.Range("K1").IndentLevel = 0

# Match the selection left behind by the edit.
$ws.Range("A2:K2").Select() | Out-Null
